$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.597.76"
$ws.Range("E2").Value = "  -3.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.897.21"
$ws.Range("E3").Value = "  -4.07%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.19"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.22"
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -2.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.895.81"
$ws.Range("E9").Value = "  -4.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.72"
$ws.Range("E10").Value = "  +5.35%  "
$ws.Range("E11").Value = "  -4.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000224"
$ws.Range("E13").Value = "  -4.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.70"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.377.78"
$ws.Range("E16").Value = "  -4.03%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.566.17"
$ws.Range("E17").Value = "  -3.07%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.80"
$ws.Range("E18").Value = "  -2.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.901.68"
$ws.Range("E19").Value = "  -3.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "424.29"
$ws.Range("E20").Value = "  -5.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.59"
$ws.Range("E21").Value = "  -4.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.670"
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("E23").Value = "  -4.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.19"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.99"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.27"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.21"
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("E32").Value = "  -3.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.39"
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0832"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.64"
$ws.Range("E37").Value = "  -3.53%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.15"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.95"
$ws.Range("E39").Value = "  -3.59%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.02"
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.125"
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.70"
$ws.Range("E42").Value = "  -3.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.292"
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.51"
$ws.Range("E44").Value = "  +3.02%  "
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "369.96"
$ws.Range("E46").Value = "  -5.31%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "133.46"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.651.27"
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.36"
$ws.Range("E50").Value = "  +6.42%  "
$ws.Range("E51").Value = "  -1.37%  "
